$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1): update column F (想去人数) values
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 6051
$ws1.Range("F9").Value = 38
$ws1.Range("F12").Value = 141
$ws1.Range("F13").Value = 343
$ws1.Range("F14").Value = 443
$ws1.Range("F15").Value = 3054
$ws1.Range("F16").Value = 4
$ws1.Range("F18").Value = 1693

# Sheet "全部类型" (sheet4): update column F (想去人数) values
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 6051
$ws4.Range("F10").Value = 38
$ws4.Range("F13").Value = 141
$ws4.Range("F14").Value = 343
$ws4.Range("F15").Value = 443
$ws4.Range("F16").Value = 3054
$ws4.Range("F17").Value = 4
$ws4.Range("F19").Value = 1693
